$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) for several event rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 105
$ws1.Range("F3").Value = 50
$ws1.Range("F6").Value = 9354
$ws1.Range("F7").Value = 842
$ws1.Range("F9").Value = 1192
$ws1.Range("F10").Value = 1121
$ws1.Range("F12").Value = 90
$ws1.Range("F14").Value = 258
$ws1.Range("F15").Value = 406
$ws1.Range("F18").Value = 1254

# Sheet "全部类型" (sheet4): same events appear here (merged with other types), mirror updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 105
$ws4.Range("F3").Value = 50
$ws4.Range("F7").Value = 9354
$ws4.Range("F8").Value = 842
$ws4.Range("F10").Value = 1192
$ws4.Range("F11").Value = 1121
$ws4.Range("F13").Value = 90
$ws4.Range("F15").Value = 258
$ws4.Range("F16").Value = 406
$ws4.Range("F19").Value = 1254
